$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename columns to short machine-readable names ---
$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'

# --- Title-case the lowercase Spanish connector words (de/del/el/los/la/y/las) ---
# --- plus two standalone casing fixes (GUANAJUATO, MonteMorelos) ---
$ws.Range('B4').Value = 'Rincón De Romos'
$ws.Range('B18').Value = 'Chiapa De Corzo'
$ws.Range('B32').Value = 'Ocozocoautla De Espinosa'
$ws.Range('B55').Value = 'Hidalgo Del Parral'
$ws.Range('B64').Value = 'San Francisco De Borja'
$ws.Range('A86').Value = 'Ciudad De México'
$ws.Range('B90').Value = 'Cuajimalpa De Morelos'
$ws.Range('B111').Value = 'Pánuco De Coronado'
$ws.Range('A119').Value = 'Estado De México'
$ws.Range('B119').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B121').Value = 'Almoloya De Juárez'
$ws.Range('B124').Value = 'Atizapán De Zaragoza'
$ws.Range('B132').Value = 'Ecatepec De Morelos'
$ws.Range('B142').Value = 'Naucalpan De Juárez'
$ws.Range('B145').Value = 'San Felipe Del Progreso'
$ws.Range('B150').Value = 'Tenango Del Valle'
$ws.Range('B152').Value = 'Tlalnepantla De Baz'
$ws.Range('B155').Value = 'Valle De Chalco Solidaridad'
$ws.Range('B156').Value = 'Villa Del Carbón'
$ws.Range('A163').Value = 'Guanajuato'
$ws.Range('B166').Value = 'Apaseo El Alto'
$ws.Range('B167').Value = 'Apaseo El Grande'
$ws.Range('B173').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B176').Value = 'Jaral Del Progreso'
$ws.Range('B185').Value = 'San Diego De La Unión'
$ws.Range('B186').Value = 'San Francisco Del Rincón'
$ws.Range('B188').Value = 'San Luis De La Paz'
$ws.Range('B189').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B191').Value = 'Silao De La Victoria'
$ws.Range('B195').Value = 'Valle De Santiago'
$ws.Range('B199').Value = 'Acapulco De Juárez'
$ws.Range('B201').Value = 'Ajuchitlán Del Progreso'
$ws.Range('B202').Value = 'Alcozauca De Guerero'
$ws.Range('B204').Value = 'Atoyac De Álvarez'
$ws.Range('B206').Value = 'Buenavista De Cuéllar'
$ws.Range('B207').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B209').Value = 'Coyuca De Benítez'
$ws.Range('B210').Value = 'Cuetzala Del Progreso'
$ws.Range('B211').Value = 'Cutzamala De Pinzón'
$ws.Range('B214').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B215').Value = 'Iguala De La Independencia'
$ws.Range('B217').Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range('B227').Value = 'Taxco De Alarcón'
$ws.Range('B229').Value = 'Técpan De Galeana'
$ws.Range('B230').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B233').Value = 'Tlapa De Comonfort'
$ws.Range('B243').Value = 'Cuautepec De Hinojosa'
$ws.Range('B248').Value = 'Huejutla De Reyes'
$ws.Range('B251').Value = 'Jacala De Ledezma'
$ws.Range('B254').Value = 'Mineral De La Reforma'
$ws.Range('B255').Value = 'Mixquiahuala De Juárez'
$ws.Range('B256').Value = 'Omitlán De Juárez'
$ws.Range('B257').Value = 'Pachuca De Soto'
$ws.Range('B259').Value = 'Santiago De Anaya'
$ws.Range('B262').Value = 'Tenango De Doria'
$ws.Range('B263').Value = 'Tepehuacán De Guerero'
$ws.Range('B264').Value = 'Tezontepec De Aldama'
$ws.Range('B268').Value = 'Tula De Allende'
$ws.Range('B269').Value = 'Tulancingo De Bravo'
$ws.Range('B270').Value = 'Zacualtipán De Ángeles'
$ws.Range('B271').Value = 'Zapotlán De Juárez'
$ws.Range('B275').Value = 'Ahualulco De Mercado'
$ws.Range('B277').Value = 'Atotonilco El Alto'
$ws.Range('B282').Value = 'Encarnación De Díaz'
$ws.Range('B284').Value = 'Ixtlahuacán De Los Membrillos'
$ws.Range('B288').Value = 'Jilotlán De Los Dolores'
$ws.Range('B292').Value = 'Lagos De Moreno'
$ws.Range('B294').Value = 'San Juan De Los Lagos'
$ws.Range('B295').Value = 'San Miguel El Alto'
$ws.Range('B297').Value = 'Tamazula De Gordiano'
$ws.Range('B300').Value = 'Tepatitlán De Morelos'
$ws.Range('B302').Value = 'Unión De San Antonio'
$ws.Range('B306').Value = 'Zacoalco De Torres'
$ws.Range('B308').Value = 'Zapotlán El Grande'
$ws.Range('B349').Value = 'Tiquicheo De Nicolás Romero'
$ws.Range('B375').Value = 'Tetela Del Volcán'
$ws.Range('B376').Value = 'Tlaltizapán De Zapata'
$ws.Range('B383').Value = 'Ixtlán Del Río'
$ws.Range('B393').Value = 'Montemorelos'
$ws.Range('B395').Value = 'San Nicolás De Los Garza'
$ws.Range('B398').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B399').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B406').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B407').Value = 'Mixistlán De La Reforma'
$ws.Range('B409').Value = 'Oaxaca De Juárez'
$ws.Range('B410').Value = 'Putla Villa De Guerero'
$ws.Range('B416').Value = 'San Dionisio Del Mar'
$ws.Range('B427').Value = 'San Pedro Y San Pablo Ayutla'
$ws.Range('B448').Value = 'Tataltepec De Valdés'
$ws.Range('B449').Value = 'Tezoatlán De Segura Y Luna'
$ws.Range('B450').Value = 'Tlacolula De Matamoros'
$ws.Range('B451').Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range('B452').Value = 'Zimatlán De Álvarez'
$ws.Range('B460').Value = 'Cuetzalan Del Progreso'
$ws.Range('B466').Value = 'Izúcar De Matamoros'
$ws.Range('B470').Value = 'Los Reyes De Juárez'
$ws.Range('B480').Value = 'Tecali De Herrera'
$ws.Range('B492').Value = 'Amealco De Bonfil'
$ws.Range('B493').Value = 'Cadereyta De Montes'
$ws.Range('B496').Value = 'Pinal De Amoles'
$ws.Range('B498').Value = 'San Juan Del Río'
$ws.Range('B503').Value = 'Axtla De Terrazas'
$ws.Range('B511').Value = 'Santa María Del Río'
$ws.Range('B538').Value = 'Jalpa De Méndez'
$ws.Range('B563').Value = 'Contla De Juan Cuamatzi'
$ws.Range('B565').Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range('B570').Value = 'San Pablo Del Monte'
$ws.Range('B582').Value = 'Cazones De Herrera'
$ws.Range('B586').Value = 'Cosamaloapan De Carpio'
$ws.Range('B587').Value = 'Cosautlán De Carvajal'
$ws.Range('B594').Value = 'Huiloapan De Cuauhtémoc'
$ws.Range('B595').Value = 'Ignacio De La Llave'
$ws.Range('B603').Value = 'Lerdo De Tejada'
$ws.Range('B604').Value = 'Martínez De La Torre'
$ws.Range('B605').Value = 'Medellín De Bravo'
$ws.Range('B608').Value = 'Paso De Ovejas'
$ws.Range('B609').Value = 'Paso Del Macho'
$ws.Range('B611').Value = 'Poza Rica De Hidalgo'
$ws.Range('B614').Value = 'Sayula De Alemán'
$ws.Range('B615').Value = 'Soledad De Doblado'
$ws.Range('B630').Value = 'Cañitas De Felipe Pescador'
$ws.Range('B632').Value = 'El Plateado De Joaquín Amaro'
$ws.Range('B642').Value = 'Moyahua De Estrada'
$ws.Range('B643').Value = 'Noria De Ángeles'

# --- Two 1-ULP floating point literal corrections ---
$ws.Range("D68").Value = 0.09128822984244672
$ws.Range("D562").Value = 0.009267840593141796

# --- Remove trailing footnote rows (656-660); dimension shrinks to A1:D654 ---
$ws.Range("A655:D660").EntireRow.Delete()

